# 2023.05.02 11:35 - Add files via upload
# Adds a new sheet "시트2" (a CREATE-TABLE style column listing for the
# T3USER/회원정보 table) after "시트1", and fixes the KIND column
# description text on 시트1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Fix the KIND description string on 시트1 (cell H39)
# ---------------------------------------------------------------------
$ws1.Range("H39").Value = "글머리 1-구매,  2-판매,  3-의뢰, 4-홍보"

# ---------------------------------------------------------------------
# 2) Insert the new sheet right after 시트1
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "시트2"

# Column A width
$ws2.Columns.Item(1).ColumnWidth = 22

# Give every cell that actually holds data the same base look
# (font/alignment) already used throughout 시트1 (style index "1" there)
# by copying the format of a cell that carries it. Rows differ in how
# many columns they use, so paste per used range rather than blanket
# filling the whole A1:D19 rectangle (that would create spurious blank
# cells the source file doesn't have).
$ws1.Range("E1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$ws2.Range("A2:C2").PasteSpecial(-4122)
$ws2.Range("A3:D7").PasteSpecial(-4122)
$ws2.Range("A8:C9").PasteSpecial(-4122)
$ws2.Range("A10:D13").PasteSpecial(-4122)
$ws2.Range("A14:C14").PasteSpecial(-4122)
$ws2.Range("A15:D19").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3) Fill in the values, row by row
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "CREATE TABLE"

$ws2.Range("A2").Value = "UNO"
$ws2.Range("B2").Value = "NUMBER"
$ws2.Range("C2").Value = "("

$ws2.Range("A3").Value = "ID"
$ws2.Range("B3").Value = "VARCHAR2"
$ws2.Range("C3").Value = "("
$ws2.Range("D3").Value = 40

$ws2.Range("A4").Value = "PASSWD"
$ws2.Range("B4").Value = "VARCHAR2"
$ws2.Range("C4").Value = "("
$ws2.Range("D4").Value = 40

$ws2.Range("A5").Value = "NAME"
$ws2.Range("B5").Value = "VARCHAR2"
$ws2.Range("C5").Value = "("
$ws2.Range("D5").Value = 20

$ws2.Range("A6").Value = "NICK"
$ws2.Range("B6").Value = "VARCHAR2"
$ws2.Range("C6").Value = "("
$ws2.Range("D6").Value = 20

$ws2.Range("A7").Value = "GENDER"
$ws2.Range("B7").Value = "CHAR"
$ws2.Range("C7").Value = "("
$ws2.Range("D7").Value = 2

$ws2.Range("A8").Value = "BIRTH"
$ws2.Range("B8").Value = "DATE"
$ws2.Range("C8").Value = "("

$ws2.Range("A9").Value = "PHONE"
$ws2.Range("B9").Value = "VARCHAR2"
$ws2.Range("C9").Value = "("

$ws2.Range("A10").Value = "EMAIL"
$ws2.Range("B10").Value = "VARCHAR2"
$ws2.Range("C10").Value = "("
$ws2.Range("D10").Value = 100

$ws2.Range("A11").Value = "ADDR"
$ws2.Range("B11").Value = "VARCHAR2"
$ws2.Range("C11").Value = "("
$ws2.Range("D11").Value = 200

$ws2.Range("A12").Value = "GRADE"
$ws2.Range("B12").Value = "NUMBER"
$ws2.Range("C12").Value = "("
$ws2.Range("D12").Value = 2

$ws2.Range("A13").Value = "CGRADE"
$ws2.Range("B13").Value = "NUMBER"
$ws2.Range("C13").Value = "("
$ws2.Range("D13").Value = 2

$ws2.Range("A14").Value = "POINT"
$ws2.Range("B14").Value = "NUMBER"
$ws2.Range("C14").Value = "("

$ws2.Range("A15").Value = "AUTH1YN"
$ws2.Range("B15").Value = "CHAR"
$ws2.Range("C15").Value = "("
$ws2.Range("D15").Value = 2

$ws2.Range("A16").Value = "AUTH2YN"
$ws2.Range("B16").Value = "CHAR"
$ws2.Range("C16").Value = "("
$ws2.Range("D16").Value = 2

$ws2.Range("A17").Value = "CREAFLG"
$ws2.Range("B17").Value = "CHAR"
$ws2.Range("C17").Value = "("
$ws2.Range("D17").Value = 2

$ws2.Range("A18").Value = "ADMINFLG"
$ws2.Range("B18").Value = "CHAR"
$ws2.Range("C18").Value = "("
$ws2.Range("D18").Value = 2

$ws2.Range("A19").Value = "STATUS"
$ws2.Range("B19").Value = "CHAR"
$ws2.Range("C19").Value = "("
$ws2.Range("D19").Value = 1

# ---------------------------------------------------------------------
# 4) Re-apply the handful of cells that use the alternate fills in the
#    source sheet (white fill "6" / light-blue fill "7"); column C keeps
#    the plain style "1" throughout, so it is never touched here.
# ---------------------------------------------------------------------
$ws1.Range("C19").Copy()
$ws2.Range("B17").PasteSpecial(-4122)

$ws1.Range("D19").Copy()
$ws2.Range("D17").PasteSpecial(-4122)
$ws2.Range("D18").PasteSpecial(-4122)

$ws1.Range("B21").Copy()
$ws2.Range("A19").PasteSpecial(-4122)

$ws1.Range("C21").Copy()
$ws2.Range("B19").PasteSpecial(-4122)

$ws1.Range("D21").Copy()
$ws2.Range("D19").PasteSpecial(-4122)

$ws1.Range("A1").Select()

Write-Output "done"
